$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "92.495.20"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.56%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.108.74"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.70%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.88"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.42%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "613.31"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.91%  "
$ws.Range("E7").Value = "  -1.81%  "
$ws.Range("E8").Value = "  -0.39%  "
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3.106.13"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.73%  "
$ws.Range("E11").Value = "  +3.65%  "
$ws.Range("E12").Value = "  -3.87%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000244"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.23%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "92.275.32"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.74%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "33.89"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.50%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.42"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.83%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.683.98"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.64%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.104.90"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.03%  "
$ws.Range("E19").Value = "  -0.20%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.52"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.61%  "
$ws.Range("E21").Value = "  -1.82%  "
$ws.Range("E22").Value = "  +0.50%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "438.59"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.35%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.13"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.59%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.20"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.16%  "
$ws.Range("E26").Value = "  -6.68%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "85.44"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.50%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.52"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.80%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.273.65"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.69%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.06%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.182"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +7.79%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.236"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.50%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.124"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -13.03%  "
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "9.16"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.20%  "
$ws.Range("B35").Value = "Binance-PegBSC-USD"
$ws.Range("C35").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -27.44%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "8.12"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +7.92%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.164"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.61%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "25.65"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.49%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.99"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.35%  "
$ws.Range("E40").Value = "  -12.26%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "23.90"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +7.72%  "
$ws.Range("E42").Value = "  -3.20%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "463.57"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.73%  "
$ws.Range("E44").Value = "  -3.31%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.31"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.36%  "
$ws.Range("E46").Value = "  +0.02%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "159.96"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.25%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.682"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.21%  "
$ws.Range("E49").Value = "  -5.16%  "
$ws.Range("E50").Value = "  -2.18%  "
$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0325"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.93%  "
